$d = $word.ActiveDocument

# 1. Arraignment sentence date: "on February 27, 2022." -> "on February 28, 2022."
$d.Content.Find.Execute(" on February 27, 2022.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " on February 28, 2022.", 2)

# 2. Sentencing table - "Fine Amount" row: "$ 45" -> "$ 4"
#    (use a precisely-bounded Range rather than Find, since Find scoped to a
#    cell Range still searches/replaces across the whole document here)
$t = $d.Tables(1)
$fineCell = $t.Cell(6, 2)
$fineStart = $fineCell.Range.Start
$fineRange = $d.Range($fineStart, $fineStart + 4)
$fineRange.Text = "$ 4"

# 3. Sentencing table - "Jail Days" row: "5" -> "2"
$jailCell = $t.Cell(8, 2)
$jailStart = $jailCell.Range.Start
$jailRange = $d.Range($jailStart, $jailStart + 1)
$jailRange.Text = "2"

# 4. "Fines and Costs paid in full by" date: "February 27, 2022" -> "February 28, 2022"
$d.Content.Find.Execute("February 27, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "February 28, 2022", 2)

# 5. Credit for days served: "3 day" -> "1 day"
$d.Content.Find.Execute("3 day", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1 day", 2)
